$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column A - labels (entered first so the shared-string table fills in the
# same order as the authored workbook)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("A2").Value = "Participated students"
$ws.Range("A3").Value = "Helpfullness of the environment"

# ---------------------------------------------------------------------------
# Column F - unit / source-of-value descriptions
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "Unit of the value"
$ws.Range("F2").Value = "Number of persons"
$ws.Range("F3").Value = "Value from the questionnaires"

# ---------------------------------------------------------------------------
# Header row (row 1) columns B:E
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Target"
$ws.Range("C1").Value = "Start"
$ws.Range("D1").Value = "End"
$ws.Range("E1").Value = "Percent"

$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Color = 0
$headerRange.Font.Bold = $true

# G1 has a style applied but no content (plain/normal black font, no bold)
$ws.Range("G1").Font.Color = 0

# ---------------------------------------------------------------------------
# Row 2 numeric values
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 80
$ws.Range("E2").Value = 1.6

# ---------------------------------------------------------------------------
# Row 3 numeric values
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = 80
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 60
$ws.Range("E3").Value = 0.75

# Data rows (2:3), columns A-F - italic font with custom teal colour
$dataRange = $ws.Range("A2:F3")
$dataRange.Font.Italic = $true
$dataRange.Font.Color = 6375440

# Percent column formatted as percentage
$ws.Range("E2:E3").NumberFormat = "0%"

# ---------------------------------------------------------------------------
# Column widths - autofit to content
# ---------------------------------------------------------------------------
$ws.Range("A1:F3").Columns.AutoFit()
